$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header "Colocação" in C1, styled like the existing header row (bold, centered)
$ws.Range("C1").Value = "Colocação"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108

# Ranking values for rows 2-8
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 25
